$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (row 59 and row 60), matching columns A:O
$newRows = @(
    @(45704.335694444446, 8, 6, 217, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3115, 4150),
    @(45707.291921296295, 8, 6, 217, 390, 366, 388, 2681, 388, 1216, 119, 304, 30, 3119, 4209)
)

$startRow = 59
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
